# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.448.58"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.872.38"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'1.019"
$ws.Range("E4").Value = "  +1.37%  "
$ws.Range("D5").Value = "'317.15"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'1.019"
$ws.Range("D7").Value = "'0.5128"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "'0.3934"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.08299"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'1.113"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'41.97"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'6.244"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.869.40"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'20.44"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.219"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.019"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'0.00001106"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "'91.23"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'0.06755"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'17.67"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.018"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "'5.971"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "28.473.45"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.265"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "2.081.89"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "'161.68"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'20.76"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("D30").Value = "'127.14"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'0.1052"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'1.033"
$ws.Range("D33").Value = "'5.800"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "'3.647"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "'0.02437"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'0.06490"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'9.160"
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'1.244"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "'1.184"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'0.6427"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "'4.988"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'11.18"
$ws.Range("D44").Value = "'0.6023"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "'12.89"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "'3.705"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "'1.220"
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("D48").Value = "'1.989"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").Value = "'121.91"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.06857"
$ws.Range("E51").Value = "  -0.70%  "
